$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "preparacao do sorteio pronta" - the new draw's columns (C:E) are filled
# in, so auto-fit them to their content the same way the other data
# columns (B, F, G, I, J) were already sized.
$ws.Columns("C:E").AutoFit() | Out-Null
$ws.Columns("C:E").ColumnWidth = 11.75

# Scroll the view down to the bottom of the table and leave the new
# column E selected, matching where the user ended up after finishing
# the data entry / formatting pass.
$ws.Range("B17").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 17
$ws.Range("E4:E25").Select() | Out-Null
